$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell 'B2' 'Bitcoin'
Set-TextCell 'C2' 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell 'D2' '28.208.93'
Set-TextCell 'E2' '  +0.56%  '
Set-TextCell 'B3' 'Ethereum'
Set-TextCell 'C3' 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell 'D3' '1.881.72'
Set-TextCell 'E3' '  +0.45%  '
Set-TextCell 'B4' 'TetherUSD'
Set-TextCell 'C4' 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 'D4' '1.002'
Set-TextCell 'E4' '  -0.18%  '
Set-TextCell 'B5' 'BNB'
Set-TextCell 'C5' 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 'D5' '313.23'
Set-TextCell 'E5' '  +0.12%  '
Set-TextCell 'B6' 'USDC'
Set-TextCell 'C6' 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 'D6' '1.002'
Set-TextCell 'E6' '  -0.12%  '
Set-TextCell 'B7' 'XRP'
Set-TextCell 'C7' 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 'D7' '0.5133'
Set-TextCell 'E7' '  +0.31%  '
Set-TextCell 'B8' 'Cardano'
Set-TextCell 'C8' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 'D8' '0.3902'
Set-TextCell 'E8' '  +2.02%  '
Set-TextCell 'B9' 'Dogecoin'
Set-TextCell 'C9' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 'D9' '0.08370'
Set-TextCell 'E9' '  +1.15%  '
Set-TextCell 'B10' 'Polygon'
Set-TextCell 'C10' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D10' '1.118'
Set-TextCell 'E10' '  +0.38%  '
Set-TextCell 'B11' 'OKB'
Set-TextCell 'C11' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D11' '41.56'
Set-TextCell 'E11' '  -0.25%  '
Set-TextCell 'B12' 'Polkadot'
Set-TextCell 'C12' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D12' '6.234'
Set-TextCell 'E12' '  +0.11%  '
Set-TextCell 'B13' 'Solana'
Set-TextCell 'C13' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 'D13' '20.75'
Set-TextCell 'E13' '  +1.17%  '
Set-TextCell 'B14' 'WrappedEther'
Set-TextCell 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D14' '1.884.17'
Set-TextCell 'E14' '  +0.08%  '
Set-TextCell 'B15' 'Chainlink'
Set-TextCell 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D15' '7.301'
Set-TextCell 'E15' '  +1.23%  '
Set-TextCell 'B16' 'BinanceUSD'
Set-TextCell 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D16' '1.001'
Set-TextCell 'E16' '  -0.28%  '
Set-TextCell 'B17' 'ShibaInu'
Set-TextCell 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D17' '0.00001107'
Set-TextCell 'E17' '  +1.09%  '
Set-TextCell 'B18' 'Litecoin'
Set-TextCell 'C18' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D18' '91.44'
Set-TextCell 'E18' '  +0.49%  '
Set-TextCell 'B19' 'TRON'
Set-TextCell 'C19' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 'D19' '0.06655'
Set-TextCell 'E19' '  +0.10%  '
Set-TextCell 'B20' 'Avalanche'
Set-TextCell 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D20' '17.75'
Set-TextCell 'E20' '  -1.29%  '
Set-TextCell 'B21' 'Dai'
Set-TextCell 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D21' '1.002'
Set-TextCell 'E21' '  -0.06%  '
Set-TextCell 'B22' 'Uniswap'
Set-TextCell 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D22' '6.052'
Set-TextCell 'E22' '  -0.02%  '
Set-TextCell 'B23' 'WrappedBTC'
Set-TextCell 'C23' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 'D23' '28.238.75'
Set-TextCell 'E23' '  +0.52%  '
Set-TextCell 'B24' 'Cosmos'
Set-TextCell 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D24' '11.19'
Set-TextCell 'E24' '  +0.28%  '
Set-TextCell 'B25' 'Toncoin'
Set-TextCell 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D25' '2.261'
Set-TextCell 'E25' '  -0.26%  '
Set-TextCell 'B26' 'WrappedliquidstakedEther2.0'
Set-TextCell 'C26' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 'D26' '2.088.45'
Set-TextCell 'E26' '  -0.46%  '
Set-TextCell 'B27' 'LidoDAOToken'
Set-TextCell 'C27' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D27' '2.512'
Set-TextCell 'E27' '  -3.18%  '
Set-TextCell 'B28' 'Monero'
Set-TextCell 'C28' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D28' '158.46'
Set-TextCell 'E28' '  +0.61%  '
Set-TextCell 'B29' 'EthereumClassic'
Set-TextCell 'C29' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D29' '20.63'
Set-TextCell 'E29' '  +0.12%  '
Set-TextCell 'B30' 'BitcoinCash'
Set-TextCell 'C30' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 'D30' '125.50'
Set-TextCell 'E30' '  -0.25%  '
Set-TextCell 'B31' 'Stellar'
Set-TextCell 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D31' '0.1065'
Set-TextCell 'E31' '  +0.77%  '
Set-TextCell 'B32' 'ImmutableX'
Set-TextCell 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D32' '1.042'
Set-TextCell 'E32' '  -0.38%  '
Set-TextCell 'B33' 'Filecoin'
Set-TextCell 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D33' '5.886'
Set-TextCell 'E33' '  +4.97%  '
Set-TextCell 'B34' 'HuobiToken'
Set-TextCell 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D34' '3.591'
Set-TextCell 'E34' '  -0.45%  '
Set-TextCell 'B35' 'FraxShare'
Set-TextCell 'C35' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D35' '9.732'
Set-TextCell 'E35' '  +0.59%  '
Set-TextCell 'B36' 'VeChain'
Set-TextCell 'C36' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D36' '0.02457'
Set-TextCell 'E36' '  +0.24%  '
Set-TextCell 'B37' 'Hedera'
Set-TextCell 'C37' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D37' '0.06556'
Set-TextCell 'E37' '  -0.27%  '
Set-TextCell 'B38' 'Algorand'
Set-TextCell 'C38' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D38' '0.2199'
Set-TextCell 'E38' '  +1.41%  '
Set-TextCell 'B39' 'ARBITRUM'
Set-TextCell 'C39' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D39' '1.211'
Set-TextCell 'E39' '  -0.20%  '
Set-TextCell 'B40' 'TheSandbox'
Set-TextCell 'C40' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D40' '0.6530'
Set-TextCell 'E40' '  +0.69%  '
Set-TextCell 'B41' 'InternetComputer(DFINITY)'
Set-TextCell 'C41' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D41' '5.018'
Set-TextCell 'E41' '  +2.69%  '
Set-TextCell 'B42' 'TrustWalletToken'
Set-TextCell 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D42' '1.231'
Set-TextCell 'E42' '  -1.31%  '
Set-TextCell 'B43' 'Aptos'
Set-TextCell 'C43' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D43' '11.31'
Set-TextCell 'E43' '  -0.30%  '
Set-TextCell 'B44' 'Decentraland'
Set-TextCell 'C44' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D44' '0.6127'
Set-TextCell 'E44' '  -0.15%  '
Set-TextCell 'B45' 'EnergySwap'
Set-TextCell 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D45' '13.13'
Set-TextCell 'E45' '  +0.51%  '
Set-TextCell 'B46' 'WEMIXTOKEN'
Set-TextCell 'C46' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D46' '1.289'
Set-TextCell 'E46' '  -0.63%  '
Set-TextCell 'B47' 'PancakeSwap'
Set-TextCell 'C47' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D47' '3.672'
Set-TextCell 'E47' '  +0.13%  '
Set-TextCell 'B48' 'NEARProtocol'
Set-TextCell 'C48' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D48' '2.019'
Set-TextCell 'E48' '  +0.33%  '
Set-TextCell 'B49' 'EOS'
Set-TextCell 'C49' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell 'D49' '1.230'
Set-TextCell 'E49' '  +0.99%  '
Set-TextCell 'B50' 'Quant'
Set-TextCell 'C50' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D50' '121.73'
Set-TextCell 'E50' '  +0.61%  '
Set-TextCell 'B51' 'Aave'
Set-TextCell 'C51' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D51' '78.23'
Set-TextCell 'E51' '  -2.95%  '
